$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing Covid risk row (row 9): rename the risk from
# "Contracting Coronavirus" to "Covid-19". The cause/effect text stay the same,
# they just end up pointing at different shared-string entries once the old
# string is dropped and the new ones are appended.
$ws.Range("C9").Value = "Covid-19"
$ws.Range("D9").Value = "This project will be carried out during a pandemic where chances of catching the virus will be extremely high"
$ws.Range("E9").Value = "The illness can severly impact my abilty to work."

# Add the new risk "Inadequate testing" in row 10.
$ws.Range("C10").Value = "Inadequate testing"
$ws.Range("D10").Value = "This may be due to insufficient unite testing or not closly adhering to the TDD methodology"
$ws.Range("E10").Value = "This will result is some aspects of the software behaving unexpectadly."

# Match styling of the other data rows: Cause/Effect columns use wrap text
# + vertical centering, same formatting already applied to D9/E9/D8/E8 etc.
$ws.Range("D10").WrapText = $true
$ws.Range("D10").VerticalAlignment = -4108
$ws.Range("E10").WrapText = $true
$ws.Range("E10").VerticalAlignment = -4108

# Row 10 grows taller to fit the new text.
$ws.Rows.Item(10).RowHeight = 64.5

# Update the active selection to match the edit (best effort - the view's
# scroll/topLeftCell position is not settable through this object model).
$ws.Range("C11").Select()
